# Refresh the scraped "cryptos" price/volume table (rows 2-51) with the
# latest values from the upstream source, matching the GitHub Actions
# update commit. Price cells (column D) get NumberFormat="@" applied
# before the assignment so Excel stores them as literal text instead of
# silently coercing/ truncating numeric-looking strings (e.g. "18.50",
# "1.00", "0.0000270") into floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.169.29"
$ws.Range("E2").Value = "  +3.79%  "
$ws.Range("D3").Value = "3.484.73"
$ws.Range("E3").Value = "  +4.83%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "555.87"
$ws.Range("E5").Value = "  +7.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.39"
$ws.Range("E6").Value = "  +7.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.651"
$ws.Range("E7").Value = "  +10.56%  "
$ws.Range("D8").Value = "3.479.90"
$ws.Range("E8").Value = "  +4.85%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.634"
$ws.Range("E10").Value = "  +5.40%  "
$ws.Range("E11").Value = "  +15.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.36"
$ws.Range("E12").Value = "  +3.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  +6.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.27"
$ws.Range("E14").Value = "  +3.81%  "
$ws.Range("D15").Value = "4.053.75"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "3.493.05"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.50"
$ws.Range("E17").Value = "  +6.42%  "
$ws.Range("E18").Value = "  +4.06%  "
$ws.Range("D19").Value = "66.236.36"
$ws.Range("E19").Value = "  +4.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.97"
$ws.Range("E20").Value = "  +7.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.990"
$ws.Range("E21").Value = "  +4.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.21"
$ws.Range("E22").Value = "  +13.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.04"
$ws.Range("E23").Value = "  +11.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.35"
$ws.Range("E24").Value = "  +6.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.12"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.92"
$ws.Range("E26").Value = "  -3.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.89"
$ws.Range("E27").Value = "  +8.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.26"
$ws.Range("E28").Value = "  +9.97%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.07"
$ws.Range("E30").Value = "  +11.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.09"
$ws.Range("E31").Value = "  +5.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "630.58"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.57"
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.70"
$ws.Range("E34").Value = "  +5.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.110"
$ws.Range("E35").Value = "  +5.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.85"
$ws.Range("E36").Value = "  +4.00%  "
$ws.Range("D37").Value = "0.0₃0806"
$ws.Range("E37").Value = "  +11.84%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("B39").Value = "InjectiveProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.63"
$ws.Range("E39").Value = "  +5.50%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.146"
$ws.Range("E40").Value = "  +18.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.384"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("E42").Value = "  +14.12%  "
$ws.Range("D43").Value = "3.107.06"
$ws.Range("E43").Value = "  +6.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.59"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.35"
$ws.Range("E46").Value = "  +12.45%  "
$ws.Range("E47").Value = "  +10.76%  "
$ws.Range("E48").Value = "  +5.68%  "
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +8.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.91"
$ws.Range("E51").Value = "  +2.57%  "
